$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.847.62'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.42%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.892.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.49%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7764'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.81%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '245.00'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.37%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3144'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.30%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.50'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.88%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07289'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.99%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08102'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.08%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7709'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.77%  '

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.504'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.29%  '

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.889.11'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.31%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.44'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.69%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.331'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.89%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.849.77'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.47%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.97'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.02%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.32'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007807'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.56%  '

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.02%  '

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.170'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.14%  '

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.142.18'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.06%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9996'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1597'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.478'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.53%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.57'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.19%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.79'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.35%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.049'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.93%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.436'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.72%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.550'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.42%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.485'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.33%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.092'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.34%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05532'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.66%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.259'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.10%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7557'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.50%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.002'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.80%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.638'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.34%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01923'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.37%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.786'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.29%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.166.84'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +13.54%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.97%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4440'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.99%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.932'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.63%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8497'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.00%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9994'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.09%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.894'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.73%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.69'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.79%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.972'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.57%  '

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.054'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.71%  '

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.486'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.99%  '
